# Apply the commit: add TIPO PERGUNTA, USUARIO and PERSONAGEM sheets with
# supporting data/formulas, and update view state (active sheet/selection).

$wb = $excel.ActiveWorkbook
$respostaSheet = $wb.Worksheets.Item(1)
$perguntaSheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Add the three new worksheets, in order, after PERGUNTA.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tipoSheet = $wb.Worksheets.Add($null, $lastSheet)
$tipoSheet.Name = "TIPO PERGUNTA"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$usuarioSheet = $wb.Worksheets.Add($null, $lastSheet)
$usuarioSheet.Name = "USUARIO"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$personagemSheet = $wb.Worksheets.Add($null, $lastSheet)
$personagemSheet.Name = "PERSONAGEM"

# ---------------------------------------------------------------------------
# 2. TIPO PERGUNTA sheet
#    String entry order matters for the shared-string table, so the B column
#    (category names) is filled before the A/B headers, matching how the
#    workbook was actually authored.
# ---------------------------------------------------------------------------
$tipoSheet.Range("B2").Value = "Banco de dados"
$tipoSheet.Range("B3").Value = "Java"
$tipoSheet.Range("B4").Value = "Programação"
$tipoSheet.Range("A1").Value = "Cod_tipo pergunta"
$tipoSheet.Range("B1").Value = "Pergunta"

$tipoSheet.Range("A2").Value = 1
$tipoSheet.Range("A3").Value = 2
$tipoSheet.Range("A4").Value = 3

$tipoSheet.Range("C2").Formula = "=IF(`$A2=`"`",`"`",`"('`"&B2&`"')`")"
$tipoSheet.Range("C3:C4").Formula = "=IF(`$A3=`"`",`"`",`"('`"&B3&`"')`")"

$tipoSheet.Columns.Item(1).AutoFit() | Out-Null
$tipoSheet.Columns.Item(2).AutoFit() | Out-Null

$tipoSheet.Range("C2:C4").Select()

# ---------------------------------------------------------------------------
# 3. USUARIO sheet
#    String entry order: headers, then A2/B2, A3/B3, B4/A4, D3, D4, F1, K1
#    (this reproduces the exact shared-string ordering of the workbook).
# ---------------------------------------------------------------------------
$usuarioSheet.Range("A1").Value = "NOME"
$usuarioSheet.Range("B1").Value = "SOBRENOME"
$usuarioSheet.Range("C1").Value = "EMAIL"
$usuarioSheet.Range("D1").Value = "APELIDO"
$usuarioSheet.Range("E1").Value = "SENHA"

$usuarioSheet.Range("A2").Value = "'Thiago"
$usuarioSheet.Range("B2").Value = "Santos"
$usuarioSheet.Range("A3").Value = "Gladson"
$usuarioSheet.Range("B3").Value = "Ameno"
$usuarioSheet.Range("B4").Value = "Campos"
$usuarioSheet.Range("A4").Value = "Amanda"
$usuarioSheet.Range("D3").Value = "Gladstone"
$usuarioSheet.Range("D4").Value = "Mandinha"
$usuarioSheet.Range("F1:J1").Merge()
$usuarioSheet.Range("F1").Value = "Strings Modelo"
$usuarioSheet.Range("F1:J1").HorizontalAlignment = -4108
$usuarioSheet.Range("K1").Value = "Inserts tabela usuario"

$usuarioSheet.Range("D2").Value = "'Thiago"

$usuarioSheet.Range("E2").Value = 123
$usuarioSheet.Range("E3").Value = 456
$usuarioSheet.Range("E4").Value = 789

# Email helper column (one formula per row; not filled as a block).
$usuarioSheet.Range("C2").Formula = "=LOWER(LEFT(A2,1)&LEFT(B2,1)&`"@gmail.com`")"
$usuarioSheet.Range("C3").Formula = "=LOWER(LEFT(A3,1)&LEFT(B3,1)&`"@gmail.com`")"
$usuarioSheet.Range("C4").Formula = "=LOWER(LEFT(A4,1)&LEFT(B4,1)&`"@gmail.com`")"

# F column - INSERT statement builder, first name fragment.
$usuarioSheet.Range("F2").Formula = "=IF(`$A2=`"`",`"`",`"('`"&A2&`"',`")"
$usuarioSheet.Range("F3:F23").Formula = "=IF(`$A3=`"`",`"`",`"('`"&A3&`"',`")"

# G/H/I columns - rows 2:3 share one formula per column.
$usuarioSheet.Range("G2:G3").Formula = "=IF(`$A2=`"`",`"`",`"'`"&B2&`"',`")"
$usuarioSheet.Range("H2:H3").Formula = "=IF(`$A2=`"`",`"`",`"'`"&C2&`"',`")"
$usuarioSheet.Range("I2:I3").Formula = "=IF(`$A2=`"`",`"`",`"'`"&D2&`"',`")"

# Row 4 - entered individually for G/H/I.
$usuarioSheet.Range("G4").Formula = "=IF(`$A4=`"`",`"`",`"'`"&B4&`"',`")"
$usuarioSheet.Range("H4").Formula = "=IF(`$A4=`"`",`"`",`"'`"&C4&`"',`")"
$usuarioSheet.Range("I4").Formula = "=IF(`$A4=`"`",`"`",`"'`"&D4&`"',`")"

# Row 5 - G:I filled together in one pass.
$usuarioSheet.Range("G5:I5").Formula = "=IF(`$A5=`"`",`"`",`"'`"&B5&`"',`")"

# Rows 6:23 - filled per column.
$usuarioSheet.Range("G6:G23").Formula = "=IF(`$A6=`"`",`"`",`"'`"&B6&`"',`")"
$usuarioSheet.Range("H6:H23").Formula = "=IF(`$A6=`"`",`"`",`"'`"&C6&`"',`")"
$usuarioSheet.Range("I6:I23").Formula = "=IF(`$A6=`"`",`"`",`"'`"&D6&`"',`")"

# J column - password fragment.
$usuarioSheet.Range("J2").Formula = "=IF(`$A2=`"`",`"`",`"'`"&E2&`"',),`")"
$usuarioSheet.Range("J3:J23").Formula = "=IF(`$A3=`"`",`"`",`"'`"&E3&`"',),`")"

# K column - concatenated INSERT statement.
$usuarioSheet.Range("K2").Formula = "=CONCATENATE(F2,G2,H2,I2,J2)"
$usuarioSheet.Range("K3:K23").Formula = "=CONCATENATE(F3,G3,H3,I3,J3)"

$usuarioSheet.Columns.Item(11).AutoFit() | Out-Null

$usuarioSheet.Range("F2").Select()

# ---------------------------------------------------------------------------
# 4. PERSONAGEM sheet stays empty.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 5. Restore view state: RESPOSTA becomes the active/selected tab, PERGUNTA's
#    selection collapses to a single cell.
# ---------------------------------------------------------------------------
$perguntaSheet.Range("C2").Select()
$respostaSheet.Activate()
$respostaSheet.Range("F2:F41").Select()
